# Rename the "M2_" sample-name prefix used throughout the accucor sample
# columns to the new "072920_XXX2_" naming convention. The rename only
# ever occurs in the header row of each sheet, but a workbook-wide
# Find & Replace mirrors how this was actually done in Excel and keeps
# every sheet (Original, Corrected, Normalized, PoolAfterDF) in sync.
$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("M2_", "072920_XXX2_")
}

# Re-fit the column widths now that several headers grew longer after the
# rename above (e.g. "M2_12_gWAT" -> "072920_XXX2_12_gWAT").
foreach ($ws in $wb.Worksheets) {
    $ws.UsedRange.EntireColumn.AutoFit()
}

# Restore a sane selection/active-cell on the sheets whose leftover
# selection no longer makes sense (one had the whole sheet selected,
# another had a stray selection from editing), and make "Original" the
# active tab again (it had drifted to "PoolAfterDF").
$wsNormalized = $wb.Worksheets.Item("Normalized")
$wsNormalized.Range("A1").Select()

$wsPoolAfterDF = $wb.Worksheets.Item("PoolAfterDF")
$wsPoolAfterDF.Range("A1").Select()

$wsOriginal = $wb.Worksheets.Item("Original")
$wsOriginal.Activate()
$wsOriginal.Range("A1").Select()
